$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 86 (this shifts the existing rows 86:196 down to 87:197,
# and naturally carries forward the existing date-number format for column D).
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new data record.
$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 45195
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100108
$ws.Cells.Item(86, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(86, 9).Value = 100108002
$ws.Cells.Item(86, 10).Value = "Mango"
$ws.Cells.Item(86, 11).Value = "Sin especificar"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 250
$ws.Cells.Item(86, 14).Value = 11000
$ws.Cells.Item(86, 15).Value = 11000
$ws.Cells.Item(86, 16).Value = 11000
$ws.Cells.Item(86, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(86, 18).Value = "Brasil"
$ws.Cells.Item(86, 19).Value = 2750
$ws.Cells.Item(86, 20).Value = 4
